# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the existing columns and filling the data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy the formatting from the neighboring "sum" header (G1)
# so the new header cell picks up the same shared style (bold, centered,
# bordered), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data cells H2:H4 - new "Save" column values, all zero.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
